$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 73
$ws.Range("A4").HorizontalAlignment = -4131
$ws.Range("B4").Value = "Set Matrix Zeroes"

$ws.Range("E8").Select()
